# This script applies the "weekly refresh" update to the 广州-漫展信息 workbook:
#  - Sheet 1 (展览):     updates "want-to-go" counts (column F) for a number of events.
#  - Sheet 2 (演出):     updates F2, removes the now-past "2024-08-24 春日计划" event row
#                         (row 3), shifting subsequent rows up.
#  - Sheet 3 (本地生活):  updates "want-to-go" counts (column F) for two events.
#  - Sheet 4 (全部类型):  mirrors sheets 1-3 (it is the concatenation of all other
#                         sheets) - updates F for rows 2-6, removes the same past event
#                         row (old row 7), and updates F for the remaining shifted rows.
#
# Because deleting a row in this environment shifts every column *except* column A
# (the manually maintained sequential index column), column A is restored to its
# correct sequential values (0, 1, 2, ...) after each row deletion.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: 展览 (Exhibition)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$sheet1Updates = @{
    2  = 48
    3  = 27070
    4  = 662
    8  = 381
    9  = 489
    10 = 199
    12 = 317
    14 = 514
    15 = 69
    16 = 1645
    17 = 268
    18 = 962
    19 = 198
    20 = 466
    21 = 12
    22 = 111
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value2 = $sheet1Updates[$row]
}

# ---------------------------------------------------------------------------
# Sheet 2: 演出 (Performance)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# Update F2 (want-to-go count for the LoveLive event) before the row shift.
$ws2.Cells.Item(2, 6).Value2 = 4529

# Row 3 (2024-08-24, 春日计划2024) is no longer upcoming - remove it; every row
# below shifts up by one.
$ws2.Rows(3).Delete()

# Restore the sequential index column (A) for the now-shifted rows 2..23.
$ws2RowCount = 23
for ($r = 2; $r -le $ws2RowCount; $r++) {
    $ws2.Cells.Item($r, 1).Value2 = $r - 1
}

# ---------------------------------------------------------------------------
# Sheet 3: 本地生活 (Local life)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Cells.Item(2, 6).Value2 = 5201
$ws3.Cells.Item(3, 6).Value2 = 280

# ---------------------------------------------------------------------------
# Sheet 4: 全部类型 (All types) - mirrors sheets 1-3 combined.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

# F-column updates for the rows that precede the deleted row (rows 2-6).
$ws4.Cells.Item(2, 6).Value2 = 48
$ws4.Cells.Item(3, 6).Value2 = 5201
$ws4.Cells.Item(4, 6).Value2 = 280
$ws4.Cells.Item(5, 6).Value2 = 27070
$ws4.Cells.Item(6, 6).Value2 = 4529

# Row 7 (2024-08-24, 春日计划2024) is removed, same as in sheet 2.
$ws4.Rows(7).Delete()

# Restore the sequential index column (A) for the now-shifted rows 2..47.
$ws4RowCount = 47
for ($r = 2; $r -le $ws4RowCount; $r++) {
    $ws4.Cells.Item($r, 1).Value2 = $r - 1
}

# F-column updates for the rows that followed the deleted row, using their
# new (post-shift) row numbers.
$sheet4PostShiftUpdates = @{
    7  = 662
    21 = 381
    22 = 489
    23 = 199
    26 = 317
    30 = 514
    31 = 69
    33 = 1645
    34 = 268
    35 = 962
    37 = 198
    38 = 466
    39 = 12
    40 = 111
}
foreach ($row in $sheet4PostShiftUpdates.Keys) {
    $ws4.Cells.Item($row, 6).Value2 = $sheet4PostShiftUpdates[$row]
}
